# Spawn entities from view
#
# Adds three new animal entries (Bird, Snake, Lion) to the "Animals" sheet
# and their matching prefab-view rows to the "PrefabsView" sheet, mirroring
# the existing "Monkey" entry/row pattern.

$wb = $excel.ActiveWorkbook

$animals = $wb.Worksheets.Item("Animals")
$prefabs = $wb.Worksheets.Item("PrefabsView")

# Style donors that already carry the correct (pre-existing) cell formats.
$animalsIdStyleSrc     = $animals.Cells.Item(2, 1)   # A2 "Monkey" -> IDS column style
$animalsNameStyleSrc   = $animals.Cells.Item(2, 2)   # B2 "Monkey" -> Name column style
$prefabsViewStyleSrc   = $animalsIdStyleSrc          # matches PrefabsView view/IDS columns (s=1)
$prefabsPathStyleSrc   = $prefabs.Cells.Item(2, 3)   # C2 prefab path column style (s=4)

$names = @("Bird", "Snake", "Lion")

# --- Animals sheet: fill IDS (A) / Name (B) columns for rows 3..5 ---
for ($i = 0; $i -lt $names.Length; $i++) {
    $name = $names[$i]
    $row = 3 + $i

    $a = $animals.Cells.Item($row, 1)
    $a.Value = $name
    $animalsIdStyleSrc.Copy()
    $a.PasteSpecial(-4122)

    $b = $animals.Cells.Item($row, 2)
    $b.Value = $name
    $animalsNameStyleSrc.Copy()
    $b.PasteSpecial(-4122)
}

# --- PrefabsView sheet: fill view name (A) / IDS (B) / prefab path (C) for rows 3..5 ---
for ($i = 0; $i -lt $names.Length; $i++) {
    $name = $names[$i]
    $row = 3 + $i

    $a = $prefabs.Cells.Item($row, 1)
    $a.Value = "$name view"
    $prefabsViewStyleSrc.Copy()
    $a.PasteSpecial(-4122)

    $b = $prefabs.Cells.Item($row, 2)
    $b.Value = $name
    $prefabsViewStyleSrc.Copy()
    $b.PasteSpecial(-4122)

    $c = $prefabs.Cells.Item($row, 3)
    $c.Value = "../Prefabs/$name.xml"
    $prefabsPathStyleSrc.Copy()
    $c.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
